$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "315.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.28%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.89%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.128"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.23%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08210"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.47%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.132"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.22%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.08%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.14%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9289"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1036"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.07%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1899"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.98%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09137"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5.08%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03619"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.28%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09908"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.51%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001432"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.42%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005756"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.60%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.467"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.08%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "7.74%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3435"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.41%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1310"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.81%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.096"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.91%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.04%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04546"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.17%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.77%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004714"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.57%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.12%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004500"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-5.35%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01966"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.08%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04920"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.88%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007653"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.46%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.26%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007862"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.54%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002130"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.99%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01181"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "4.62%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006746"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.25%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.25%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "38.52"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-22.10%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-15.08%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.25%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.25%"
